# Generate Report for Handback
# Updates timestamps / status values across the "Overview", "zh-cn" and "de-de"
# sheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview.Range("G4").Value = "2016-09-04 04:18:57"
$wsOverview.Range("G5").Value = "2016-09-04 04:18:57"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime (H)
$wsZhCn.Range("H4").Value = "2016-09-04 04:18:52"
$wsZhCn.Range("H5").Value = "2016-09-04 04:18:52"

# Correspond Handback DateTime (K)
$wsZhCn.Range("K4").Value = "2016-09-04 04:19:16"
$wsZhCn.Range("K5").Value = "2016-09-04 04:19:16"

# --- de-de sheet ---
# Priority column (E): ht -> mt
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# Correspond Handback DateTime (K)
$wsDeDe.Range("K4").Value = "2016-09-04 04:19:22"
$wsDeDe.Range("K5").Value = "2016-09-04 04:19:22"
